$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$th = $m.Theme
Write-Output "Theme type: $($th.GetType())"
try { Write-Output "Theme.Name: $($th.Name)" } catch { Write-Output "no Name" }
try { Write-Output "Theme.Parent: $($th.Parent)" } catch { Write-Output "no Parent" }
